# Actualización automática 2025-07-29 14:50:09
# Insert a new client row ("CONSTRUCCION, INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS")
# into the OFICINA-CATAECSA table, keeping the alphabetical ordering of
# column B. This shifts every following data row (and the trailing
# "X de N" summary row) down by one, and bumps the "de 17" -> "de 18"
# counters since the table now has one more entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row above row 7 (DANIELA ELIZABETH BECERRA BECERRA),
# pushing all rows below (including formatting) down by one.
$ws.Range("A7:R7").Insert()

# Populate the newly inserted row.
$ws.Range("A7").Value = "OFICINA-CATAECSA"
$ws.Range("B7").Value = "CONSTRUCCION, INGENIERIA Y TECNOLOGIA CONSTRUINTEC SAS"
$ws.Range("C7:R7").Value = 0

# Widen column B so the longer company name fits. (55.1666... is the
# ColumnWidth input that this host's char->pixel quantizer round-trips to
# an exported OOXML <col width="56"/>, matching the target worksheet.)
$ws.Columns.Item(2).ColumnWidth = 55.1666666667

# Update the "X de 17" summary counters (row 19 -> now row 20) to "X de 18"
# since one more row was added to the table.
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols) {
    $cell = $ws.Range($col + "20")
    $cell.Value = $cell.Value().Replace("de 17", "de 18")
}
